$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Expand the table from 2 columns (Goals, Progress) to 4 columns ---
$table = $ws.ListObjects.Item("Table1")
$table.Resize($ws.Range("A1:D24"))

# Rename header cells through the table's header row range so the
# ListColumn names (and the underlying table XML) are updated together
# with the worksheet header cells.
$hdr = $table.HeaderRowRange
$hdr.Cells.Item(1, 1).Value = "Goals"
$hdr.Cells.Item(1, 2).Value = "Progress in code"
$hdr.Cells.Item(1, 3).Value = "Progress in report"
$hdr.Cells.Item(1, 4).Value = "Progress in video"

# Row 2: URLs - move "In Progress" from column B to column C
$ws.Range("B2").Value = $null
$ws.Range("C2").Value = "In Progress"

# Row 3: Source-code - change status from "In Progress" to "Done" (still col B)
$ws.Range("B3").Value = "Done"

# Row 4: Demonstration video - add "In Progress" in column D
$ws.Range("D4").Value = "In Progress"

# Row 6: Problem statement - move "Done" out of B, set "In Progress" in C
$ws.Range("B6").Value = $null
$ws.Range("C6").Value = "In Progress"

# Row 7: Existing approaches - move "Done" out of B, set "In Progress" in C
$ws.Range("B7").Value = $null
$ws.Range("C7").Value = "In Progress"

# Row 8: Similarities and differences - move "Done" out of B, set "In Progress" in C
$ws.Range("B8").Value = $null
$ws.Range("C8").Value = "In Progress"

# Rows 10 and 11 (Data analysis / Data pre-processing) keep their "Done" in col B - unchanged

# Row 12: Applied machine learning algorithms - add "Done" in column B
$ws.Range("B12").Value = "Done"

# Row 13: Model tuning - add "Done" in column B
$ws.Range("B13").Value = "Done"

# Row 14: Evaluation - add "Done" in column B
$ws.Range("B14").Value = "Done"

# Row 16: Analysis and evaluation - add "Done" in column B
$ws.Range("B16").Value = "Done"

# Row 17: Conclusion - add "Done" in column B
$ws.Range("B17").Value = "Done"

# Rows 19-22: move "In progress" from column B to column C
$ws.Range("B19").Value = $null
$ws.Range("C19").Value = "In progress"

$ws.Range("B20").Value = $null
$ws.Range("C20").Value = "In progress"

$ws.Range("B21").Value = $null
$ws.Range("C21").Value = "In progress"

$ws.Range("B22").Value = $null
$ws.Range("C22").Value = "In progress"

# Row 23 (Complete source code as text at Appendix B) - unchanged, no second column

# --- Column widths for the new/expanded columns ---
# (input values chosen so the engine's internal rounding lands as close as
# possible to the target stored widths of 16.81640625 / 18.1796875 / 17.36328125)
$ws.Columns.Item(2).ColumnWidth = 16.0
$ws.Columns.Item(3).ColumnWidth = 17.333333333333332
$ws.Columns.Item(4).ColumnWidth = 16.5

# --- Selection moves to D4 ---
$ws.Range("D4").Select()
